# Updated symbol list on Tue Jan  3 15:40:18 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values to the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "245.40";      E = "-0.56%" }
    @{ Row = 3;  D = "29.13";       E = "-1.23%" }
    @{ Row = 4;  D = "5.256";       E = "1.13%" }
    @{ Row = 5;  D = "0.05714";     E = "0.06%" }
    @{ Row = 6;  D = "6.616";       E = "0.55%" }
    @{ Row = 7;  D = "3.179";       E = "3.85%" }
    @{ Row = 8;  D = "0.8550";      E = "-0.55%" }
    @{ Row = 9;  D = "0.8577";      E = "-2.28%" }
    @{ Row = 10;                    E = "0.50%" }
    @{ Row = 11; D = "0.07079";     E = "0.22%" }
    @{ Row = 12; D = "0.03177";     E = "10.76%" }
    @{ Row = 13; D = "0.09290";     E = "-1.09%" }
    @{ Row = 14; D = "0.001532";    E = "1.28%" }
    @{ Row = 15; D = "0.0005945";   E = "-94.24%" }
    @{ Row = 16; D = "0.006082";    E = "-0.69%" }
    @{ Row = 17; D = "3.511";       E = "0.84%" }
    @{ Row = 18; D = "2.181";       E = "-0.29%" }
    @{ Row = 19; D = "0.3171";      E = "-0.16%" }
    @{ Row = 20;                    E = "1.56%" }
    @{ Row = 22; D = "3.483";       E = "0.44%" }
    @{ Row = 23; D = "0.04121";     E = "-0.89%" }
    @{ Row = 24;                    E = "-3.56%" }
    @{ Row = 25;                    E = "0.32%" }
    @{ Row = 26; D = "0.004144";    E = "-18.03%" }
    @{ Row = 27;                    E = "-0.65%" }
    @{ Row = 28;                    E = "-25.19%" }
    @{ Row = 40; D = "0.03809";     E = "1.48%" }
    @{ Row = 41; D = "0.1063";      E = "-1.02%" }
    @{ Row = 42; D = "0.002419";    E = "-4.69%" }
    @{ Row = 43;                    E = "-15.05%" }
    @{ Row = 44; D = "0.009416";    E = "-5.28%" }
    @{ Row = 45; D = "0.00005291";  E = "-1.65%" }
    @{ Row = 46; D = "0.00000000751"; E = "0.17%" }
    @{ Row = 47; D = "0.08997";     E = "26.82%" }
    @{ Row = 48; D = "0.002453";    E = "-5.27%" }
    @{ Row = 49; D = "0.00002102";  E = "0.17%" }
    @{ Row = 50; D = "0.0002002";   E = "0.17%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
